$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert "update pathways part 1": restore the prior 11-row data set
# (rows 2-12), replacing the current 12-row data set (rows 2-13).
$data = @(
    @("PALM",     "58262248", "Cat",    "Feature Room 2",           "4/8/2025"),
    @("COLLIE",   "57884999", "Cat",    "If The Fur Fits",          "2/20/2025"),
    @("SPECK",    "52249653", "Dog",    "Dog Adoptions D",          "4/16/2025"),
    @("LAUREL",   "58289985", "Cat",    "Cat Adoption Condo Rooms", "4/11/2025"),
    @("Gyarados", "58096306", "Cat",    "Cat Treatment",            "3/28/2025"),
    @("Mochi",    "58353916", "Cat",    "Offsite Adoptions",        "4/22/2025"),
    @("SUDS",     "58598619", "Dog",    "Dog Adoptions A",          "5/29/2025"),
    @("Beau",     "58677023", "Dog",    "Dog Adoptions C",          "6/10/2025"),
    @("RAYNE",    "57710656", "Dog",    "If The Fur Fits",          "6/13/2025"),
    @("HEATH",    "58654173", "Dog",    "If The Fur Fits",          "6/6/2025"),
    @("MARINA",   "58706705", "Rabbit", "Adoptions Lobby",          "6/13/2025")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]

    # Columns A (Name), C (Species), D (Location) are plain text - assign
    # directly.
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]

    # Columns B (AID) and E (Intake Date) look numeric/date-like, but in
    # the workbook they are stored as literal text (shared strings), not
    # as a real number/date. Build each as a text formula, then collapse
    # it down to its literal value via copy/paste-special so the stored
    # cell stays plain text instead of becoming a number or date serial.
    $bcell = $ws.Cells.Item($row, 2)
    $bcell.Formula = '="' + $vals[1] + '"'
    $bcell.Copy()
    $bcell.PasteSpecial(-4163)

    $ecell = $ws.Cells.Item($row, 5)
    $ecell.Formula = '="' + $vals[4] + '"'
    $ecell.Copy()
    $ecell.PasteSpecial(-4163)
}

# The reverted table has one fewer row than the pre-revert version, so
# drop the now-extra trailing row.
$ws.Range("A13:E13").Delete()
